$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.393.52"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "2.667.04"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.22"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "2.666.30"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.76"
$ws.Range("E14").Value = "  -4.12%  "
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "67.460.95"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "2.682.65"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.71"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.79"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.69"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  -5.47%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.70"
$ws.Range("E26").Value = "  -4.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000102"
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "547.80"
$ws.Range("E31").Value = "  -9.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("E33").Value = "  -5.21%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -6.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.41"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.48"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("E41").Value = "  -4.91%  "
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.91"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -7.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.28"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("E47").Value = "  -5.81%  "
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.35"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.85"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("E51").Value = "  -5.03%  "
